$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "before"
$ws.Range("B1").Value = "after"

# Update the active selection to match the target workbook state.
$ws.Range("B7").Select()
